# Update Financials: Net Receivables (row 43) and Inventory (row 44)
# on the DAKT sheet's Balance Sheet section (columns D:J).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("DAKT")

# Row 43 - Net Receivables
$ws.Range("D43").Value = 84500
$ws.Range("E43").Value = 81700
$ws.Range("F43").Value = 85500
$ws.Range("G43").Value = 90200
$ws.Range("H43").Value = 90200
$ws.Range("I43").Value = 68100
$ws.Range("J43").Value = 78700

# Row 44 - Inventory
$ws.Range("D44").Value = 106300
$ws.Range("E44").Value = 102900
$ws.Range("F44").Value = 100000
$ws.Range("G44").Value = 99500
$ws.Range("H44").Value = 95600
$ws.Range("I44").Value = 88400
$ws.Range("J44").Value = 77900
